# Ajustes a las mediciones de performances segun ultimos cambios
#
# Updates the raw "Ms Totales" (column D) and "ms espera / msg" (column G)
# measurements on the "NodoSocket" sheet for the three blocks of rows
# affected by the latest benchmark run: rows 21-23 (MensajeVortex en socket -
# memory-only numbers block just above) and rows 29-31 (MensajeVortex en
# socket block). Columns E and F are formulas (msg/s and ms/msg) and
# recalculate automatically once D changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NodoSocket")

# --- Rows 21-23 -----------------------------------------------------------
$ws.Range("D21").Value = 4222.4483840000003
$ws.Range("G21").Value = 10.117475170000001

$ws.Range("D22").Value = 452.89974000000001
$ws.Range("G22").Value = 9.5499764860000003

$ws.Range("D23").Value = 98.912851000000003
$ws.Range("G23").Value = 7.3004350889999996

# --- Rows 29-31 -------------------------------------------------------------
$ws.Range("D29").Value = 12137.010151
$ws.Range("G29").Value = 1834.5638976170601

$ws.Range("D30").Value = 1583.176633
$ws.Range("G30").Value = 268.70276339669999

$ws.Range("D31").Value = 394.179149
$ws.Range("G31").Value = 181.454856963

# --- Refresh view/selection state ------------------------------------------
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C57").Select() | Out-Null
